# Updates the "广州-漫展信息" workbook to the newer scrape snapshot:
#  - drop the expired "广州·一人之下快闪店" listing (was the oldest row) from every
#    sheet that still carries it, shifting the remaining rows up
#  - refresh the running "want to go" counter (column F) for listings that are
#    still open, on every sheet that carries that listing
#
# Column A in every sheet is just the zero-based row rank (0 for the header,
# 1, 2, 3, ... for the data rows) - it is NOT part of the scraped record, so
# after a row is removed it must be renumbered back to a contiguous sequence
# rather than left shifted.

$wb = $excel.ActiveWorkbook

function Renumber-ColumnA($ws) {
    $lastRow = $ws.Cells.Item(1, 1).End(4).Row
    for ($r = 2; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }
}

function Apply-FUpdates($ws, $updates) {
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibitions) - no expired row present, only counter refreshes
# ---------------------------------------------------------------------------
$wsExhibit = $wb.Worksheets.Item(1)
$exhibitUpdates = @{
    2  = 27
    3  = 1813
    4  = 1813
    6  = 1479
    7  = 850
    8  = 377
    9  = 732
    10 = 13146
    11 = 13021
    12 = 988
    13 = 762
    16 = 64
    17 = 628
    18 = 2057
    19 = 51
    20 = 26
    21 = 36
    23 = 185
    24 = 271
    25 = 733
}
Apply-FUpdates $wsExhibit $exhibitUpdates

# ---------------------------------------------------------------------------
# Sheet "演出" (Performances) - no expired row present, only counter refresh
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item(2)
$showUpdates = @{
    7 = 99
}
Apply-FUpdates $wsShow $showUpdates

# ---------------------------------------------------------------------------
# Sheet "本地生活" (Local life) - drop the expired row, renumber column A
# ---------------------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item(3)
$wsLocal.Rows.Item(2).Delete()
Renumber-ColumnA $wsLocal

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All types) - drop the expired row, renumber column A,
# then refresh the counters (row numbers below are post-deletion positions)
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Rows.Item(2).Delete()
Renumber-ColumnA $wsAll

$allUpdates = @{
    3  = 27
    4  = 1813
    5  = 1813
    7  = 1479
    8  = 850
    9  = 377
    11 = 732
    12 = 13146
    13 = 13021
    14 = 988
    15 = 762
    18 = 64
    19 = 628
    22 = 2057
    23 = 51
    24 = 26
    25 = 36
    29 = 185
    30 = 271
    31 = 733
    32 = 99
}
Apply-FUpdates $wsAll $allUpdates
